# --------------------------------------------------------------------
# TeachingNotes.docx edit: reword the "read-modify-write" sentence,
# flip Normal style's overflowPunct off, and mint the ListLabel273..308
# character styles that a LibreOffice round-trip regenerated.
# --------------------------------------------------------------------

$d = $word.ActiveDocument

# 1) Wording change in the atomic-ops paragraph -----------------------
$d.Content.Find.Execute(
    "that the read-modify-write routine is the", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "to the read-modify-write routine and it being the", 2) | Out-Null

# 2) Normal style: stop letting punctuation overflow the margin -------
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.HangingPunctuation = 0

# 3) Re-mint the ListLabel character styles 273-308 --------------------
$fonts = @(
    "Symbol","Courier New","Wingdings","Wingdings","Wingdings","Wingdings",
    "Wingdings","Wingdings","Wingdings",
    "Symbol","Courier New","Wingdings","Wingdings","Wingdings","Wingdings",
    "Wingdings","Wingdings","Wingdings",
    "OpenSymbol","OpenSymbol","OpenSymbol","OpenSymbol","OpenSymbol","OpenSymbol",
    "OpenSymbol","OpenSymbol","OpenSymbol",
    "OpenSymbol","OpenSymbol","OpenSymbol","OpenSymbol","OpenSymbol","OpenSymbol",
    "OpenSymbol","OpenSymbol","OpenSymbol"
)

for ($i = 0; $i -lt $fonts.Length; $i++) {
    $num = 273 + $i
    $font = $fonts[$i]
    $style = $d.Styles.Add("ListLabel $num", 2)
    $style.QuickStyle = $true
    $style.Font.NameBi = $font
    if ($font -ne "OpenSymbol") {
        $style.Font.Size = 10
    }
}
